$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 held the text "R40" (row label for the 4th rule). The commit
# replaces its content with the text "1", while keeping its existing
# style/format untouched.
#
# A plain `Range.Value = "1"` would be auto-typed by Excel as the *number*
# 1 (since the cell's number format is General), which also pulls in a
# new/derived style. To keep the literal string "1" (matching the source
# workbook's shared-string cell type) with the original style, write it as
# a formula that evaluates to the text "1", then convert that formula to a
# static value in place via copy / paste-values.
$cell = $ws.Range("B11")
$cell.Formula = '=TEXT(1,"0")'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
